$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1666666666666667
$ws.Range("C2").Value = 0.6145833333333334
$ws.Range("J2").Value = 0.03125
$ws.Range("P2").Value = 0.09722222222222222
$ws.Range("S2").Value = 0.09027777777777778
$ws.Range("B3").Value = 0.01098901098901099
$ws.Range("C3").Value = 0.02197802197802198
$ws.Range("J3").Value = 0.03296703296703297
$ws.Range("P3").Value = 0.8076923076923077
$ws.Range("S3").Value = 0.1263736263736264
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.673469387755102
$ws.Range("S4").Value = 0.3061224489795918
$ws.Range("B6").Value = 0.06172839506172839
$ws.Range("D6").Value = 0.01646090534979424
$ws.Range("E6").Value = 0.00411522633744856
$ws.Range("F6").Value = 0.06584362139917696
$ws.Range("J6").Value = 0.242798353909465
$ws.Range("O6").Value = 0.01646090534979424
$ws.Range("Q6").Value = 0.1769547325102881
$ws.Range("R6").Value = 0.05761316872427984
$ws.Range("S6").Value = 0.3580246913580247
$ws.Range("B7").Value = 0.1117021276595745
$ws.Range("D7").Value = 0.02127659574468085
$ws.Range("F7").Value = 0.02659574468085106
$ws.Range("J7").Value = 0.1329787234042553
$ws.Range("O7").Value = 0.03723404255319149
$ws.Range("Q7").Value = 0.2234042553191489
$ws.Range("R7").Value = 0.101063829787234
$ws.Range("S7").Value = 0.3457446808510639
$ws.Range("B8").Value = 0.0821256038647343
$ws.Range("D8").Value = 0.02093397745571659
$ws.Range("E8").Value = 0.001610305958132045
$ws.Range("F8").Value = 0.06924315619967794
$ws.Range("J8").Value = 0.1288244766505636
$ws.Range("O8").Value = 0.01127214170692432
$ws.Range("Q8").Value = 0.1658615136876006
$ws.Range("R8").Value = 0.107890499194847
$ws.Range("S8").Value = 0.4122383252818035
$ws.Range("B9").Value = 0.08968609865470852
$ws.Range("D9").Value = 0.0179372197309417
$ws.Range("F9").Value = 0.05381165919282511
$ws.Range("J9").Value = 0.05381165919282511
$ws.Range("O9").Value = 0.008968609865470852
$ws.Range("Q9").Value = 0.2062780269058296
$ws.Range("R9").Value = 0.1390134529147982
$ws.Range("S9").Value = 0.4304932735426009
$ws.Range("B10").Value = 0.09491778774289986
$ws.Range("D10").Value = 0.02017937219730942
$ws.Range("E10").Value = 0.0007473841554559044
$ws.Range("F10").Value = 0.07100149476831091
$ws.Range("J10").Value = 0.1136023916292975
$ws.Range("O10").Value = 0.02466367713004484
$ws.Range("Q10").Value = 0.2010463378176383
$ws.Range("R10").Value = 0.09865470852017937
$ws.Range("S10").Value = 0.375186846038864
$ws.Range("G11").Value = 0.1137123745819398
$ws.Range("J11").Value = 0.1003344481605351
$ws.Range("K11").Value = 0.1638795986622074
$ws.Range("L11").Value = 0.6086956521739131
$ws.Range("S11").Value = 0.01337792642140468
$ws.Range("G12").Value = 0.7021276595744681
$ws.Range("J12").Value = 0.2287234042553191
$ws.Range("L12").Value = 0.0425531914893617
$ws.Range("S12").Value = 0.02659574468085106
$ws.Range("G13").Value = 0.7435897435897436
$ws.Range("J13").Value = 0.2051282051282051
$ws.Range("S13").Value = 0.05128205128205128
$ws.Range("F15").Value = 0.01639344262295082
$ws.Range("H15").Value = 0.2540983606557377
$ws.Range("I15").Value = 0.0778688524590164
$ws.Range("J15").Value = 0.3360655737704918
$ws.Range("K15").Value = 0.06967213114754098
$ws.Range("O15").Value = 0.02868852459016394
$ws.Range("S15").Value = 0.2172131147540984
$ws.Range("F16").Value = 0.01477832512315271
$ws.Range("H16").Value = 0.2019704433497537
$ws.Range("I16").Value = 0.08374384236453201
$ws.Range("J16").Value = 0.4729064039408867
$ws.Range("K16").Value = 0.09359605911330049
$ws.Range("M16").Value = 0.004926108374384237
$ws.Range("O16").Value = 0.03448275862068965
$ws.Range("S16").Value = 0.09359605911330049
$ws.Range("F17").Value = 0.01984126984126984
$ws.Range("H17").Value = 0.2341269841269841
$ws.Range("I17").Value = 0.08531746031746032
$ws.Range("J17").Value = 0.4047619047619048
$ws.Range("K17").Value = 0.07539682539682539
$ws.Range("M17").Value = 0.0119047619047619
$ws.Range("O17").Value = 0.0615079365079365
$ws.Range("S17").Value = 0.1071428571428571
$ws.Range("F18").Value = 0.01140684410646388
$ws.Range("H18").Value = 0.2357414448669201
$ws.Range("I18").Value = 0.09125475285171103
$ws.Range("J18").Value = 0.3422053231939163
$ws.Range("K18").Value = 0.1064638783269962
$ws.Range("M18").Value = 0.01140684410646388
$ws.Range("O18").Value = 0.09505703422053231
$ws.Range("S18").Value = 0.1064638783269962
$ws.Range("F19").Value = 0.01804511278195489
$ws.Range("H19").Value = 0.2443609022556391
$ws.Range("I19").Value = 0.08947368421052632
$ws.Range("J19").Value = 0.3609022556390977
$ws.Range("K19").Value = 0.1075187969924812
$ws.Range("M19").Value = 0.02180451127819549
$ws.Range("N19").Value = 0.001503759398496241
$ws.Range("O19").Value = 0.07368421052631578
$ws.Range("S19").Value = 0.08270676691729323
